# Update workbook with correct forecast output
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B ("ASIN" etc. shift right by one)
$ws1.Columns.Item(2).Insert()

# Header row
$ws1.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week labels (strip leading zero for single-digit weeks) and week start dates
$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekDates  = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

# Keep the new date column as plain text (not auto-converted to a date serial)
$ws1.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $weekLabels[$i]
    $ws1.Cells.Item($row, 2).Value = $weekDates[$i]
}

# Updated MyForecast values (now column D after the insert)
$myForecast = @{ 2 = 124; 5 = 132; 9 = 126; 11 = 144; 15 = 121 }
foreach ($row in $myForecast.Keys) {
    $ws1.Cells.Item($row, 4).Value = $myForecast[$row]
}

# is_holiday_week column (now column J) should hold boolean values instead of numbers
for ($row = 2; $row -le 17; $row++) {
    $ws1.Cells.Item($row, 10).Value = $false
}

# ---------------------------------------------------------------------
# Sheet 2: "Summary"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

# These cells hold numeric-looking text; keep them as text (not auto-converted to numbers)
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B14").NumberFormat = "@"

$ws2.Range("B9").Value = "2119"
$ws2.Range("B10").Value = "1052"
$ws2.Range("B14").Value = "121"
